# Apply the edits described by the diff to the active workbook/sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared strings (B1 / B5 header labels) ---
$ws.Range("B1").Value = "Tình yêu"
$ws.Range("B5").Value = "Tình bạn"

# --- Numeric data updates in column A / B ---
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = 2.1
$ws.Range("B4").Value = 3

$ws.Range("A6").Value = 3.1

$ws.Range("A7").Value = 3.2
$ws.Range("B7").Value = 2

$ws.Range("A8").Value = 4.1
$ws.Range("B8").Value = 3

# --- Column width (A:A, closest achievable increment to 11.8117647058824) ---
$ws.Range("A:A").ColumnWidth = 10.95

# --- Selection / active cell moves from B1 to B5 ---
[void]$ws.Range("B5").Select()

# --- Window tab ratio (198 -> 395, raw OOXML units; COM TabRatio is raw/330) ---
$win = $excel.ActiveWindow
$win.TabRatio = 395 / 330
